$wb = $excel.ActiveWorkbook

# The "93237bcb-bb20-4cdf-8d92-8e150c11553e.md" file has moved from
# "Ready for handoff" to "In Translation" status. Update the per-locale
# status sheets (row 2 = that file) as well as the Overview summary sheet.

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Status column is column B on the locale sheets (zh-cn / de-de); row 2
# corresponds to file 93237bcb-bb20-4cdf-8d92-8e150c11553e.md
$wsZhCn.Range("B2").Value = "In Translation"
$wsDeDe.Range("B2").Value = "In Translation"

# Overview sheet mirrors the status per locale in columns B (zh-cn) and
# C (de-de); row 2 is the same file.
$wsOverview.Range("B2").Value = "In Translation"
$wsOverview.Range("C2").Value = "In Translation"
